# EPBDS-14346 Fix conversion when the argument type is String
# Inserts a new "Method String str2str(String data)" block into the sheet,
# right before the existing "Datatype Complex" block, following the same
# visual pattern used by the other method blocks on the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new rows at row 32, pushing the "Datatype Complex" block
# (currently rows 32-36) down to rows 36-40.
$ws.Rows("32:35").Insert()

# The header/body cells of a method block are merged across B:C - merge
# before copying the formatting over so the paste brings back the exact
# border styling used by the other method blocks.
$ws.Range("B33:C33").Merge()
$ws.Range("B34:C34").Merge()

# Copy the formatting of an existing method block ("oneArg", rows 28-31)
# into the freshly inserted rows so the new block matches the visual
# style (borders/fills) used throughout the sheet.
$ws.Range("A28:E31").Copy()
$ws.Range("A32:E35").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Make sure the row heights match the rest of the sheet.
$ws.Rows("32:35").RowHeight = 13.55

# Fill in the new method block content.
$ws.Range("B33").Value = "Method String str2str(String data)"
$ws.Range("B34").Value = "return data;"

$wb.Save()
